# Deploying to gh-pages — refresh the X-LANCE member roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fix a batch of mis-tagged "state" / "degree" cells that were wrong
#    in the previous export (许洪深-style placeholders corrected to the
#    real values for the people who already left / whose degree track
#    is joint Undergrad+Master / Undergrad+PhD).
# ---------------------------------------------------------------------
$ws.Range("F97").Value  = "离开"
$ws.Range("F145").Value = "离开"

$ws.Range("D188").Value = "UP"
$ws.Range("D189").Value = "UM"
$ws.Range("D190").Value = "UM"
$ws.Range("F192").Value = "离开"

$ws.Range("D214").Value = "UM"
$ws.Range("D218").Value = "UM"
$ws.Range("D219").Value = "UM"
$ws.Range("D222").Value = "UM"

$ws.Range("F227").Value = "离开"
$ws.Range("F238").Value = "离开"
$ws.Range("F239").Value = "离开"

# ---------------------------------------------------------------------
# 2) Append the new member row (274) that fell off the roster.
# ---------------------------------------------------------------------
$ws.Range("A274").Value = "马文杰"
$ws.Range("B274").Value = "Wenjie Ma"
$ws.Range("D274").Value = "U"
$ws.Range("E274").Value = "/assets/img/members/student/马文杰.jpg"
$ws.Range("F274").Value = "离开"

# ---------------------------------------------------------------------
# 3) Re-stretch the autofilter / filter-database over the now-larger
#    A1:F274 table (it used to only cover the header row).
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:F274").AutoFilter()

$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "=Sheet1!`$A`$1:`$F`$274"

# ---------------------------------------------------------------------
# 4) Leave the selection where the new last row is, like the author did.
# ---------------------------------------------------------------------
[void]$ws.Range("E274").Select()
